$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column U: "DF_chg" header plus per-province values (rows 2-35).
# Most provinces default to 0 (formatted like the neighboring FTT/IF_val/IF_chg
# columns with the "0.00" number format / style 18). A few provinces (East Java
# row 10, Jakarta row 14, South Sulawesi row 27) carry real computed values;
# rows 14 and 27 keep the default (unstyled) format, matching row 7/10/14/27/30
# in column T which also have no explicit style.
$ws.Range("U1").Value = "DF_chg"

$specialValues = @{
    10 = -65.550286084840735
    14 = -66.488825953857457
    27 = -54.007451575054375
}
$unstyledRows = @(14, 27)

for ($r = 2; $r -le 35; $r++) {
    if ($specialValues.ContainsKey($r)) {
        $value = $specialValues[$r]
    } else {
        $value = 0
    }

    $cell = $ws.Range("U$r")
    $cell.Value = $value

    if ($unstyledRows -notcontains $r) {
        $cell.NumberFormat = "0.00"
    }
}

# Selection moved to W6 as part of this edit.
$ws.Range("W6").Select()
